$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3297984293646418
$ws.Range("C2").Value = 0.4599280677128398
$ws.Range("D2").Value = 0.3388381866556515
$ws.Range("E2").Value = 0.5820980902353584
$ws.Range("F2").Value = 0.4977639660768899

$ws.Range("B3").Value = 0.3132943648870622
$ws.Range("C3").Value = 0.3132943648870622
$ws.Range("D3").Value = 0.2339416816796298
$ws.Range("E3").Value = 0.4836751819967919
$ws.Range("F3").Value = 0.3884274887423822

$ws.Range("B4").Value = 0.3533454843703288
$ws.Range("C4").Value = 0.3533454843703288
$ws.Range("D4").Value = 0.3385089747862772
$ws.Range("E4").Value = 0.5818152411086163
$ws.Range("F4").Value = 0.5063468496531306

$ws.Range("B5").Value = 0.2173850498843069
$ws.Range("C5").Value = 0.2173850498843069
$ws.Range("D5").Value = 0.07346767206689039
$ws.Range("E5").Value = 0.2710492059883046
$ws.Range("F5").Value = 0.2289603116423796
